# Refresh the cryptocurrency snapshot: update the latest Price (column D)
# and 1-hour Volume change percentage (column E) for each affected coin row.
# A leading apostrophe forces Excel to keep the numeric-looking price strings
# as text (matching the workbook's existing inline-string cell format), and
# resetting the style back to Normal avoids leaving a stray text number format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.853.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.37%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.533.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +3.98%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  -0.08%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'568.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.77%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'147.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +4.94%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  -0.01%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'  -0.60%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'2.533.19"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +3.92%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'  +0.25%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'  -1.87%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'  +0.85%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "'  +1.45%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'27.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +5.76%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'2.986.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +3.94%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'62.861.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.15%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = "'  -0.55%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'2.536.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +4.31%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = "'  +2.93%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'336.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.94%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Value = "'  +1.60%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'6.77"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.64%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = "'  -0.09%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'65.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.09%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = "'  -3.42%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'1.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +2.61%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'  +13.50%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  +0.08%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "'  +1.43%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'7.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +10.24%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'0.0₃0811"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.81%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'1.85"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.91%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'177.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.91%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "'  +7.28%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'412.27"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +9.83%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'0.399"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.15%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'18.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.16%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'4.41"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.14%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Value = "'  +0.00%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = "'  +1.32%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'  +0.03%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'39.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.64%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'152.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +4.86%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'3.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.82%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'20.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.90%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Value = "'  +1.94%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.0966"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.49%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "'  +0.35%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = "'  +5.87%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'18.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +3.22%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Value = "'  +2.18%  "
$ws.Range("E51").Style = "Normal"

